$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates: force text to avoid Excel auto-numeric conversion
# which would corrupt multi-dot values and drop trailing zeros / formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.751.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.455.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.900.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.601.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.453.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "645.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0958"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.605"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0906"
$ws.Range("D51").Style = "Normal"

# Column E (Volume/1h) updates: plain text percentage strings
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  -5.64%  "
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -3.07%  "
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("E51").Value = "  -1.30%  "
